# Updated cryptos list on Sat Dec 16 14:09:04 UTC 2023 with GitHub Actions
# Refreshes the live price/volume snapshot pulled from coinranking.com;
# rows 48/49 and 51 also pick up new coin entries for that rank.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = '42.469.34'
$ws.Range("E2").Value = '  +0.30%  '

# Row 3 - Ethereum
$ws.Range("D3").Value = '2.246.69'
$ws.Range("E3").Value = '  -0.15%  '

# Row 4 - TetherUSD
$ws.Range("E4").Value = '  +0.15%  '

# Row 5 - BNB
$ws.Range("D5").Value = '''246.26'
$ws.Range("E5").Value = '  -0.51%  '

# Row 6 - XRP
$ws.Range("D6").Value = '''0.631'
$ws.Range("E6").Value = '  -0.29%  '

# Row 7 - Solana
$ws.Range("D7").Value = '''75.86'
$ws.Range("E7").Value = '  -0.89%  '

# Row 8 - USDC
$ws.Range("E8").Value = '  +0.13%  '

# Row 9 - Cardano
$ws.Range("D9").Value = '''0.621'
$ws.Range("E9").Value = '  -2.45%  '

# Row 10 - Avalanche
$ws.Range("D10").Value = '''44.00'
$ws.Range("E10").Value = '  +8.45%  '

# Row 11 - Dogecoin
$ws.Range("D11").Value = '''0.0949'
$ws.Range("E11").Value = '  -0.56%  '

# Row 12 - Polkadot
$ws.Range("D12").Value = '''7.29'
$ws.Range("E12").Value = '  +0.27%  '

# Row 13 - TRON
$ws.Range("E13").Value = '  -1.26%  '

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = '2.582.47'
$ws.Range("E14").Value = '  -0.16%  '

# Row 15 - Chainlink
$ws.Range("D15").Value = '''14.60'
$ws.Range("E15").Value = '  -2.44%  '

# Row 16 - Polygon
$ws.Range("D16").Value = '''0.855'
$ws.Range("E16").Value = '  -0.72%  '

# Row 17 - WrappedEther
$ws.Range("D17").Value = '2.229.22'
$ws.Range("E17").Value = '  -1.67%  '

# Row 18 - WrappedBTC
$ws.Range("D18").Value = '42.254.19'
$ws.Range("E18").Value = '  +0.02%  '

# Row 19 - ShibaInu
$ws.Range("E19").Value = '  +3.67%  '

# Row 20 - Uniswap
$ws.Range("E20").Value = '  +0.16%  '

# Row 21 - Litecoin
$ws.Range("D21").Value = '''72.23'
$ws.Range("E21").Value = '  +0.65%  '

# Row 22 - ImmutableX
$ws.Range("E22").Value = '  +1.79%  '

# Row 23 - BitcoinCash
$ws.Range("D23").Value = '''231.79'
$ws.Range("E23").Value = '  -0.11%  '

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").Value = '''9.16'
$ws.Range("E24").Value = '  +30.63%  '

# Row 25 - Dai
$ws.Range("E25").Value = '  +0.00%  '

# Row 26 - Cosmos
$ws.Range("D26").Value = '''11.45'
$ws.Range("E26").Value = '  +2.48%  '

# Row 27 - WEMIXToken
$ws.Range("E27").Value = '  -2.96%  '

# Row 28 - PancakeSwap
$ws.Range("D28").Value = '''2.31'
$ws.Range("E28").Value = '  -0.67%  '

# Row 29 - Toncoin
$ws.Range("E29").Value = '  +1.14%  '

# Row 30 - Monero
$ws.Range("D30").Value = '''168.09'
$ws.Range("E30").Value = '  -0.02%  '

# Row 31 - EthereumClassic
$ws.Range("D31").Value = '''20.69'
$ws.Range("E31").Value = '  +0.27%  '

# Row 32 - Hedera
$ws.Range("D32").Value = '''0.0827'
$ws.Range("E32").Value = '  -3.18%  '

# Row 33 - Kaspa
$ws.Range("E33").Value = '  -0.68%  '

# Row 34 - InjectiveProtocol
$ws.Range("D34").Value = '''30.55'
$ws.Range("E34").Value = '  -4.43%  '

# Row 35 - Filecoin
$ws.Range("D35").Value = '''5.31'
$ws.Range("E35").Value = '  +10.90%  '

# Row 36 - Stellar
$ws.Range("E36").Value = '  -0.65%  '

# Row 37 - RenderToken
$ws.Range("E37").Value = '  +1.14%  '

# Row 38 - VeChain
$ws.Range("D38").Value = '''0.0318'
$ws.Range("E38").Value = '  +6.42%  '

# Row 39 - Celestia
$ws.Range("D39").Value = '''13.78'
$ws.Range("E39").Value = '  +5.49%  '

# Row 40 - LidoDAOToken
$ws.Range("E40").Value = '  -1.88%  '

# Row 41 - THORChain
$ws.Range("D41").Value = '''5.80'
$ws.Range("E41").Value = '  -3.08%  '

# Row 42 - MultiversX
$ws.Range("D42").Value = '''63.67'
$ws.Range("E42").Value = '  +5.50%  '

# Row 43 - Algorand
$ws.Range("E43").Value = '  -0.71%  '

# Row 44 - Aave
$ws.Range("D44").Value = '''107.99'
$ws.Range("E44").Value = '  -7.81%  '

# Row 45 - FraxShare
$ws.Range("D45").Value = '''8.80'
$ws.Range("E45").Value = '  +0.67%  '

# Row 46 - Cronos
$ws.Range("E46").Value = '  +1.20%  '

# Row 47 - BinanceUSD
$ws.Range("E47").Value = '  -0.02%  '

# Row 48 - TrustWalletToken -> ARBITRUM (row re-ranked)
$ws.Range("B48").Value = 'ARBITRUM'
$ws.Range("C48").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D48").Value = '''1.13'
$ws.Range("E48").Value = '  -0.33%  '

# Row 49 - ARBITRUM -> TrustWalletToken (row re-ranked)
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").Value = '''1.19'
$ws.Range("E49").Value = '  +0.47%  '

# Row 50 - NEARProtocol
$ws.Range("D50").Value = '''2.36'
$ws.Range("E50").Value = '  +5.57%  '

# Row 51 - WOONetwork -> SynthetixNetwork (row re-ranked)
$ws.Range("B51").Value = 'SynthetixNetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range("D51").Value = '''4.12'
$ws.Range("E51").Value = '  -0.93%  '
